$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sending = @("ECs","ECs","FAPs","FAPs","sCs","sCs")
$target  = @("FAPs","sCs","FAPs","sCs","FAPs","sCs")
$values = @(
    @(2,0.6666666666666666,0.9305633333333333,2.79169,0.01768777137856805,0.01768777137856806,3,1,0.092904,0.278712,0.9572435868814848,0.9572435868814848,0.08645305591999999,0.77807750328,0.01693150571836015,0.01693150571836015),
    @(2,0.6666666666666666,0.9305633333333333,2.79169,0.01768777137856805,0.01768777137856806,1,0.3333333333333333,0.004149666666666667,0.012449,0.04275641311851519,0.04275641311851518,0.003861527645555556,0.03475374881,0.0007562656602079046,0.0007562656602079046),
    @(3,1,44.154177,132.462531,0.8392647337471152,0.8392647337471153,3,1,0.092904,0.278712,0.9572435868814848,0.9572435868814848,4.102099660008,36.918896940072,0.8033807840752228,0.803380784075223),
    @(3,1,44.154177,132.462531,0.8392647337471152,0.8392647337471153,1,0.3333333333333333,0.004149666666666667,0.012449,0.04275641311851519,0.04275641311851518,0.183225116491,1.649026048419,0.03588394967189232,0.03588394967189231),
    @(3,1,7.525807,22.577421,0.1430474948743168,0.1430474948743168,3,1,0.092904,0.278712,0.9572435868814848,0.9572435868814848,0.699177573528,6.292598161752001,0.1369312970879018,0.1369312970879018),
    @(3,1,7.525807,22.577421,0.1430474948743168,0.1430474948743168,1,0.3333333333333333,0.004149666666666667,0.012449,0.04275641311851519,0.04275641311851518,0.03122959044766667,0.281066314029,0.006116197786414973,0.006116197786414972)
)

for ($i = 0; $i -lt 6; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $sending[$i]
    $ws.Cells.Item($r, 2).Value = "Ntn1"
    $ws.Cells.Item($r, 3).Value = "Dcc"
    $ws.Cells.Item($r, 4).Value = $target[$i]
    $rowVals = $values[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $ws.Cells.Item($r, 5 + $j).Value = $rowVals[$j]
    }
}
